# Adds two new weekly-survey columns (15.-21.11.2021 and 22.-28.11.2021)
# to both worksheets, and refreshes the 'aktualizace' date in the footer rows.

$wb = $excel.ActiveWorkbook

# ---- Sheet 'data': add columns BR (col 70) and BS (col 71) ----
$wsData = $wb.Worksheets.Item("data")

# Header row 1: copy the format of the previous header cell (BQ1) so the new
# header cells keep the bold/border/center style (s="1"), then set the text.
$wsData.Range("BQ1").Copy()
$wsData.Range("BR1:BS1").PasteSpecial(-4122)
$wsData.Range("BR1").Value = '15.–21. 11. 2021'
$wsData.Range("BS1").Value = '22.–28. 11. 2021'

# Data rows 2-77 for column BR (col 70)
$dataBR = @(0.78,0.05,0.09,0.08,0.77,0.04,0.11,0.08,0.46,0.03,0.44,0.07000000000000001,0.85,0.07000000000000001,0.02,0.06,0.87,0.03,0.01,0.09,0.89,0.08,0.005,0.025,0.82,0.05,0.07000000000000001,0.06,0.8,0.04,0.07000000000000001,0.09,0.65,0.04,0.18,0.13,0.6899999999999999,0.08,0.13,0.1,0.74,0.07000000000000001,0.14,0.05,0.83,0.03,0.07000000000000001,0.07000000000000001,0.68,0.08,0.1,0.14,0.76,0.05,0.08,0.11,0.79,0.05,0.08,0.08,0.79,0.05,0.12,0.04,0.84,0.02,0.05,0.09,0.76,0.06,0.1,0.08,0.71,0.08,0.16,0.05)
# Data rows 2-77 for column BS (col 71)
$dataBS = @(0.76,0.06,0.1,0.08,0.76,0.05,0.1,0.09,0.42,0.09,0.45,0.04,0.85,0.05,0.04,0.06,0.87,0.05,0.015,0.065,0.82,0.1,0.005,0.075,0.8100000000000001,0.05,0.07000000000000001,0.07000000000000001,0.79,0.05,0.08,0.08,0.64,0.05,0.17,0.14,0.6899999999999999,0.09,0.13,0.09,0.73,0.08,0.15,0.04,0.8100000000000001,0.04,0.07000000000000001,0.08,0.6899999999999999,0.08,0.09,0.14,0.73,0.07000000000000001,0.07000000000000001,0.13,0.8,0.06,0.1,0.04,0.78,0.04,0.13,0.05,0.85,0.02,0.04,0.09,0.74,0.05,0.12,0.09,0.68,0.12,0.15,0.05)
for ($i = 0; $i -lt $dataBR.Length; $i++) {
    $r = $i + 2
    $wsData.Cells.Item($r, 70).Value = $dataBR[$i]
    $wsData.Cells.Item($r, 71).Value = $dataBS[$i]
}

# Footer row 78 label: bump the 'aktualizace' date
$wsData.Range("A78").Value = 'Život během pandemie, Home office, % respondentů celkově a ve skupinách, aktualizace 8. 12. 2021'

# ---- Sheet 'pocetR': add columns BQ (col 69) and BR (col 70) ----
$wsPocet = $wb.Worksheets.Item("pocetR")

$wsPocet.Range("BP1").Copy()
$wsPocet.Range("BQ1:BR1").PasteSpecial(-4122)
$wsPocet.Range("BQ1").Value = '15.–21. 11. 2021'
$wsPocet.Range("BR1").Value = '22.–28. 11. 2021'

# Data rows 2-20 for column BQ (col 69)
$pocetBQ = @(1018,286,102,250,143,91,494,249,127,148,273,560,185,387,315,316,374,400,244)
# Data rows 2-20 for column BR (col 70)
$pocetBR = @(1018,286,102,250,143,91,494,249,127,148,273,560,185,387,315,316,374,400,244)
for ($i = 0; $i -lt $pocetBQ.Length; $i++) {
    $r = $i + 2
    $wsPocet.Cells.Item($r, 69).Value = $pocetBQ[$i]
    $wsPocet.Cells.Item($r, 70).Value = $pocetBR[$i]
}

# Footer row 21: bump the label. (BQ21/BR21 stay blank, same as the rest of
# this placeholder row - Excel has no way to persist a literal empty-string
# cell distinct from "no cell" via COM Value assignment, so we simply leave
# them unset; the used range already reaches BR21 because of the data rows
# above, so the sheet dimension comes out the same either way.)
$wsPocet.Range("A21").Value = 'Život během pandemie, Home office, velikost dotázaného souboru celkově a ve skupinách, aktualizace 8. 12. 2021'

